# Refresh crypto price/volume figures in the "cryptos" sheet (Price = column D,
# Volume(1h) = column E), mirroring the upstream GitHub Actions data-refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.033.39'
$ws.Range("E2").Value = '  +1.31%  '
$ws.Range("D3").Value = '3.385.25'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''571.62'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").Value = '''141.66'
$ws.Range("E6").Value = '  +0.86%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +2.20%  '
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("D11").Value = '''0.388'
$ws.Range("E11").Value = '  -1.40%  '
$ws.Range("D12").Value = '3.963.24'
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("E13").Value = '  +2.17%  '
$ws.Range("D14").Value = '''27.91'
$ws.Range("E14").Value = '  -0.76%  '
$ws.Range("D15").Value = '3.384.33'
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("D16").Value = '''0.0000171'
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").Value = '61.119.37'
$ws.Range("E17").Value = '  +1.08%  '
$ws.Range("E18").Value = '  -2.25%  '
$ws.Range("D19").Value = '''13.65'
$ws.Range("E19").Value = '  -3.06%  '
$ws.Range("D20").Value = '''8.95'
$ws.Range("E20").Value = '  -2.78%  '
$ws.Range("D21").Value = '''384.69'
$ws.Range("E21").Value = '  -0.80%  '
$ws.Range("D22").Value = '''75.64'
$ws.Range("E22").Value = '  +2.99%  '
$ws.Range("D23").Value = '''0.553'
$ws.Range("E23").Value = '  -1.39%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("E25").Value = '  -0.55%  '
$ws.Range("D26").Value = '3.521.26'
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  +3.48%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").Value = '''7.22'
$ws.Range("E29").Value = '  -2.06%  '
$ws.Range("E30").Value = '  -0.88%  '
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("E33").Value = '  -3.80%  '
$ws.Range("D34").Value = '''23.21'
$ws.Range("E34").Value = '  -2.05%  '
$ws.Range("D35").Value = '''6.96'
$ws.Range("E35").Value = '  +0.60%  '
$ws.Range("D36").Value = '''166.00'
$ws.Range("E36").Value = '  -0.74%  '
$ws.Range("D37").Value = '3.418.97'
$ws.Range("E37").Value = '  +0.34%  '
$ws.Range("D38").Value = '''4.98'
$ws.Range("E38").Value = '  +1.52%  '
$ws.Range("E39").Value = '  -2.10%  '
$ws.Range("D40").Value = '''0.0767'
$ws.Range("E40").Value = '  -0.94%  '
$ws.Range("D41").Value = '''26.70'
$ws.Range("E41").Value = '  -1.17%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("E44").Value = '  -2.01%  '
$ws.Range("E45").Value = '  -2.15%  '
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("D47").Value = '2.456.24'
$ws.Range("E47").Value = '  -2.86%  '
$ws.Range("E48").Value = '  -0.43%  '
$ws.Range("E49").Value = '  -2.54%  '
$ws.Range("D50").Value = '''2.13'
$ws.Range("E50").Value = '  +10.73%  '
$ws.Range("E51").Value = '  -1.17%  '
